# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" stats table:
#   - Several countries' rows are re-labelled because the underlying
#     country list was re-sorted (Catar/Irlanda, Barein/Kazajistan,
#     Senegal/Islandia/Cuba/Estonia, Sierra Leona/Mauricio/Isla de
#     Man/Montenegro/Republica del Chad/Benin, Belice/Nueva Caledonia).
#   - Updated case/death/recovery figures are entered for the affected
#     rows (Espana, India, Catar, Irlanda, Kuwait, Barein, Kazajistan,
#     Senegal, Islandia, Cuba, Estonia, Sierra Leona, Mauricio, Isla de
#     Man, Montenegro, Republica del Chad, Benin, Belice, Nueva
#     Caledonia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Range("A$row").Value = $pais
    $ws.Range("B$row").Value = $casosTotales
    $ws.Range("C$row").Value = $nuevosCasos
    $ws.Range("D$row").Value = $casosActivos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $casosCriticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Espana
Set-Row 5 "España" 268143 3480 177846 63553 1650 123 26744

# India
Set-Row 16 "India" 67714 553 21150 44349 0 3 2215

# Catar / Irlanda swap position 30/31
Set-Row 30 "Catar" 23623 1103 2840 20769 72 0 14
Set-Row 31 "Irlanda" 22996 0 17110 4428 72 0 1458

# Kuwait
Set-Row 49 "Kuwait" 9286 598 2907 6314 131 7 65

# Barein / Kazajistan swap position 59/60
Set-Row 59 "Barein" 5157 216 2152 2997 2 0 8
Set-Row 60 "Kazajistan" 5138 48 1941 3166 31 0 31

# Senegal moves ahead of Islandia, Cuba, Estonia (rows 81-84 shift down one)
Set-Row 81 "Senegal" 1886 177 715 1152 6 0 19
Set-Row 82 "Islandia" 1801 0 1773 18 0 0 10
Set-Row 83 "Cuba" 1766 0 1193 496 5 0 77
Set-Row 84 "Estonia" 1741 2 751 929 5 1 61

# Sierra Leona moves ahead of Mauricio, Isla de Man, Montenegro, Republica
# del Chad, Benin (rows 130-135 shift down one)
Set-Row 130 "Sierra Leona" 338 31 72 247 0 1 19
Set-Row 131 "Mauricio" 332 0 320 2 0 0 10
Set-Row 132 "Isla de Man" 330 0 271 36 21 0 23
Set-Row 133 "Montenegro" 324 0 290 25 2 0 9
Set-Row 134 "Republica del Chad" 322 0 53 238 0 0 31
Set-Row 135 "Benin" 319 0 62 255 0 0 2

# Belice / Nueva Caledonia swap position 192/193
Set-Row 192 "Belice" 18 0 16 0 0 0 2
Set-Row 193 "Nueva Caledonia" 18 0 18 0 0 0 0
